$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chamada")
$r = $ws.Range("A1:B1")
$r.HorizontalAlignment = -4131  # left
$r.HorizontalAlignment = -4108  # center
